$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'40.996.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.42%  "
$ws.Range("D3").Value = "'2.154.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.89%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'235.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").Value = "'0.603"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.68%  "
$ws.Range("D7").Value = "'68.77"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.58%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "'0.561"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.89%  "
$ws.Range("D10").Value = "'38.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.73%  "
$ws.Range("D11").Value = "'0.0900"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.66%  "
$ws.Range("D12").Value = "'54.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.85%  "
$ws.Range("E13").Value = "  -2.94%  "
$ws.Range("D14").Value = "'6.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.21%  "
$ws.Range("D15").Value = "'2.478.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.75%  "
$ws.Range("D16").Value = "'14.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").Value = "'2.162.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.94%  "
$ws.Range("D18").Value = "'0.777"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.34%  "
$ws.Range("D19").Value = "'40.810.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("D20").Value = "'0.0₃0980"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.64%  "
$ws.Range("D21").Value = "'69.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.07%  "
$ws.Range("D22").Value = "'5.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.95%  "
$ws.Range("D23").Value = "'223.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("D24").Value = "'9.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -14.70%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'1.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -10.45%  "
$ws.Range("D27").Value = "'10.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.24%  "
$ws.Range("D28").Value = "'3.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.43%  "
$ws.Range("D29").Value = "'2.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.64%  "
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").Value = "'167.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "'19.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.18%  "
$ws.Range("D33").Value = "'29.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").Value = "'0.0745"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.29%  "
$ws.Range("D35").Value = "'5.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -11.49%  "
$ws.Range("E36").Value = "  -4.72%  "
$ws.Range("D37").Value = "'0.100"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.35%  "
$ws.Range("D38").Value = "'4.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.82%  "
$ws.Range("D39").Value = "'0.0276"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.91%  "
$ws.Range("D40").Value = "'2.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.47%  "
$ws.Range("D41").Value = "'11.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -16.66%  "
$ws.Range("D42").Value = "'5.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.28%  "
$ws.Range("D43").Value = "'57.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -12.25%  "
$ws.Range("D44").Value = "'0.185"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.47%  "
$ws.Range("D45").Value = "'8.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.30%  "
$ws.Range("D46").Value = "'0.0953"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.95%  "
$ws.Range("D47").Value = "'95.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.06%  "
$ws.Range("D48").Value = "'1.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.28%  "
$ws.Range("E49").Value = "  -5.19%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "'2.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.24%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'2.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.56%  "
